# JackieK-WorkLog.xlsx — "Add files via upload"
#
# Adds a new work-log entry (row 29: 2025-03-03, 4 hours, description of
# finishing the scanning commands / python scripts) below the existing
# last row (row 28), and normalizes row 28's "Number of Hours" cell back
# onto the same (non-filled) style used by every other row in that column
# — it had accidentally picked up a one-off "apply fill" variant of the
# style. Also nudges the sheet/window scroll position and selection to
# where the author last left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Fix up B28's style -----------------------------------------------
# B28 ("Number of Hours" for the last existing row) was carrying a
# redundant cell-format variant that only differed from the style used by
# every other row in that column (B2:B27) by an explicit (no-op, none)
# fill override. Clearing the interior pattern collapses it back onto the
# shared style.
$ws.Range("B28").Interior.Pattern = -4142

# --- Append the new row --------------------------------------------------
# Copy the fully-formatted previous row down (so font/border/alignment/
# number-format match exactly), then overwrite with the new row's values.
$ws.Range("A28:C28").Copy()
$ws.Range("A29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A29").Value = 45719
$ws.Range("B29").Value = 4
$ws.Range("C29").Value = "Finalized the commands which will be used for scanning and developed python scripts"

# --- View state ------------------------------------------------------
# Scroll so row 16 is near the top and leave the selection where the
# author last clicked (just past the new data, on the description column).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("C31").Select()

Write-Output "done"
